$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

# Update grade threshold values (minLevel/maxLevel columns)
$ws.Range("D7").Value = 250
$ws.Range("C8").Value = 251
$ws.Range("D8").Value = 350
$ws.Range("C9").Value = 351
$ws.Range("D9").Value = 500
$ws.Range("C10").Value = 501
$ws.Range("D10").Value = 800
$ws.Range("C11").Value = 801

# Update the selected cell / active selection on the sheet
$ws.Range("C10").Select()
